$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price (D) cells that will hold numeric-looking text so Excel
# keeps them as text instead of auto-converting to numbers.
$textCells = @("D4","D5","D6","D7","D9","D10","D11","D13","D14","D16","D18","D19","D20","D21","D22","D25","D26","D28","D29","D30","D31","D32","D33","D35","D36","D37","D38","D40","D41","D43","D45","D46","D47","D49","D50","D51")
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "26.463.59"
$ws.Range("E2").Value = "  +2.29%  "

# Row 3
$ws.Range("D3").Value = "1.670.82"
$ws.Range("E3").Value = "  +1.92%  "

# Row 4
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "219.33"
$ws.Range("E5").Value = "  +2.72%  "

# Row 6
$ws.Range("D6").Value = "0.5268"
$ws.Range("E6").Value = "  +1.37%  "

# Row 7
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("E8").Value = "  +2.85%  "

# Row 9
$ws.Range("D9").Value = "0.06376"
$ws.Range("E9").Value = "  +0.96%  "

# Row 10
$ws.Range("D10").Value = "21.70"
$ws.Range("E10").Value = "  +5.66%  "

# Row 11
$ws.Range("D11").Value = "0.07799"
$ws.Range("E11").Value = "  +1.66%  "

# Row 12
$ws.Range("D12").Value = "1.725.45"
$ws.Range("E12").Value = "  +5.14%  "

# Row 13
$ws.Range("D13").Value = "4.467"
$ws.Range("E13").Value = "  +1.34%  "

# Row 14
$ws.Range("D14").Value = "0.5534"
$ws.Range("E14").Value = "  +1.25%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8298"
$ws.Range("E15").Value = "  +1.87%  "

# Row 16
$ws.Range("D16").Value = "65.47"
$ws.Range("E16").Value = "  +1.92%  "

# Row 17
$ws.Range("D17").Value = "26.475.89"
$ws.Range("E17").Value = "  +2.25%  "

# Row 18
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.06%  "

# Row 19
$ws.Range("D19").Value = "4.740"
$ws.Range("E19").Value = "  +1.50%  "

# Row 20
$ws.Range("D20").Value = "192.81"
$ws.Range("E20").Value = "  +2.52%  "

# Row 21
$ws.Range("D21").Value = "10.33"
$ws.Range("E21").Value = "  +2.36%  "

# Row 22
$ws.Range("D22").Value = "6.271"
$ws.Range("E22").Value = "  +0.80%  "

# Row 23
$ws.Range("E23").Value = "  +0.12%  "

# Row 24
$ws.Range("E24").Value = "  +2.13%  "

# Row 25
$ws.Range("D25").Value = "138.13"
$ws.Range("E25").Value = "  -3.35%  "

# Row 26
$ws.Range("D26").Value = "7.393"
$ws.Range("E26").Value = "  +0.79%  "

# Row 27
$ws.Range("E27").Value = "  +2.53%  "

# Row 28
$ws.Range("D28").Value = "1.419"
$ws.Range("E28").Value = "  +1.38%  "

# Row 29
$ws.Range("D29").Value = "0.06167"
$ws.Range("E29").Value = "  +5.11%  "

# Row 30
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  +1.99%  "

# Row 31
$ws.Range("D31").Value = "3.611"
$ws.Range("E31").Value = "  +6.82%  "

# Row 32
$ws.Range("D32").Value = "3.393"
$ws.Range("E32").Value = "  +0.43%  "

# Row 33
$ws.Range("D33").Value = "1.679"
$ws.Range("E33").Value = "  +2.92%  "

# Row 34
$ws.Range("E34").Value = "  +1.94%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.6056"
$ws.Range("E35").Value = "  +8.65%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.419"
$ws.Range("E36").Value = "  +1.04%  "

# Row 37
$ws.Range("D37").Value = "2.770"
$ws.Range("E37").Value = "  +1.05%  "

# Row 38
$ws.Range("D38").Value = "0.01609"
$ws.Range("E38").Value = "  +1.06%  "

# Row 39
$ws.Range("D39").Value = "1.095.06"
$ws.Range("E39").Value = "  +7.72%  "

# Row 40
$ws.Range("D40").Value = "6.028"
$ws.Range("E40").Value = "  +3.85%  "

# Row 41
$ws.Range("D41").Value = "0.8549"
$ws.Range("E41").Value = "  +0.46%  "

# Row 42
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("D43").Value = "100.65"
$ws.Range("E43").Value = "  +2.19%  "

# Row 44
$ws.Range("D44").Value = "1.814.22"
$ws.Range("E44").Value = "  +1.41%  "

# Row 45
$ws.Range("D45").Value = "58.03"
$ws.Range("E45").Value = "  +4.86%  "

# Row 46
$ws.Range("D46").Value = "8.172"
$ws.Range("E46").Value = "  +2.27%  "

# Row 47
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  -4.61%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.05204"
$ws.Range("E49").Value = "  +1.27%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.481"
$ws.Range("E50").Value = "  +8.75%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.4233"
$ws.Range("E51").Value = "  +0.50%  "

# Restore the default (Normal) style on the cells we temporarily reformatted
# as text, so no stray number-format style lingers on them.
foreach ($addr in $textCells) {
  $ws.Range($addr).Style = "Normal"
}

